# Apply the "expand age structure (M = 0.065)" update:
#   - add 4 new scenario rows (20-23) to the OM sheet, cloned from the
#     existing template row 19 (same style / constants) with the
#     scenario-specific fields overwritten
#   - move the "active" sheet/tab selection from EM_2Fl to OM (matches
#     the tabSelected + selection changes in the diff)

$wb = $excel.ActiveWorkbook
$om = $wb.Worksheets.Item(1)     # "OM" sheet

# ---------------------------------------------------------------------
# 1) Build out rows 20-23 on the OM sheet.
#    Start from row 19, which already carries the correct shared
#    formatting/style (font on I/J, general number formats, etc.) and
#    the columns that stay constant across every scenario row
#    (H, I, J, K, L, M, N, O, P).
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row=20; A="Fast_LG_O_High_Ext";     B=100; C=30; D=2; E="Const_Ramp_Const"; F="Const_Ramp_Const"; G=0.065;
       Q="logistic"; R="gamma";   S=3.5;  T=0.65; U=5.5; V=0.75; W=15.5; X=8;    Y=19; Z=8;   AA="F_Vary" },
    @{ Row=21; A="Fast_LG_Y_High_Ext";     B=100; C=30; D=2; E="Const_Ramp_Const"; F="Const_Ramp_Const"; G=0.065;
       Q="logistic"; R="gamma";   S=3.5;  T=0.65; U=5.5; V=0.75; W=5;    X=5;    Y=7;  Z=6.5; AA="F_Vary" },
    @{ Row=22; A="Fast_GL_O_High_Rev_Ext"; B=100; C=30; D=2; E="Const_Ramp_Const"; F="Const_Ramp_Const"; G=0.065;
       Q="gamma";    R="logistic"; S=15.5; T=8;    U=19;  V=8;    W=3.5;  X=0.65; Y=5.5; Z=0.75; AA="F_Vary" },
    @{ Row=23; A="Fast_GL_Y_High_Rev_Ext"; B=100; C=30; D=2; E="Const_Ramp_Const"; F="Const_Ramp_Const"; G=0.065;
       Q="gamma";    R="logistic"; S=5;    T=5;    U=7;   V=6.5;  W=3.5;  X=0.65; Y=5.5; Z=0.75; AA="F_Vary" }
)

foreach ($rowData in $newRows) {
    $targetRow = $rowData.Row

    # Clone formatting + the constant columns from the template row.
    $src = $om.Range("A19:AA19")
    $dst = $om.Range("A" + $targetRow + ":AA" + $targetRow)
    $src.Copy($dst)

    # Overwrite the scenario-specific cells.
    $om.Range("A$targetRow").Value2 = $rowData.A
    $om.Range("B$targetRow").Value2 = $rowData.B
    $om.Range("C$targetRow").Value2 = $rowData.C
    $om.Range("D$targetRow").Value2 = $rowData.D
    $om.Range("E$targetRow").Value2 = $rowData.E
    $om.Range("F$targetRow").Value2 = $rowData.F
    $om.Range("G$targetRow").Value2 = $rowData.G
    $om.Range("Q$targetRow").Value2 = $rowData.Q
    $om.Range("R$targetRow").Value2 = $rowData.R
    $om.Range("S$targetRow").Value2 = $rowData.S
    $om.Range("T$targetRow").Value2 = $rowData.T
    $om.Range("U$targetRow").Value2 = $rowData.U
    $om.Range("V$targetRow").Value2 = $rowData.V
    $om.Range("W$targetRow").Value2 = $rowData.W
    $om.Range("X$targetRow").Value2 = $rowData.X
    $om.Range("Y$targetRow").Value2 = $rowData.Y
    $om.Range("Z$targetRow").Value2 = $rowData.Z
    $om.Range("AA$targetRow").Value2 = $rowData.AA

    $om.Rows.Item($targetRow).RowHeight = 20
}

# ---------------------------------------------------------------------
# 2) Update the active sheet/selection to match the saved view in the
#    diff: the OM sheet becomes the active/selected tab (cursor on
#    G23), which automatically clears the old tabSelected flag that
#    used to sit on EM_2Fl (its own selection, A11, is untouched).
# ---------------------------------------------------------------------
$om.Range("G23").Select() | Out-Null
$om.Activate()
